$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "שם העסק"
$ws.Range("B1").Value = "אמצעי זיהוי התשלום"
$ws.Range("C1").Value = "תאריך התשלום"
$ws.Range("D1").Value = "תאריך החיוב בחשבון"
$ws.Range("E1").Value = "סכום"

# Row 2
$ws.Range("A2").Value = "רמי לוי"
$ws.Range("B2").Value = 1234
$ws.Range("B2").NumberFormat = "@"
$ws.Range("C2").Value = 45306
$ws.Range("C2").NumberFormat = "m/d/yy"
$ws.Range("D2").Value = 45306
$ws.Range("D2").NumberFormat = "m/d/yy"
$ws.Range("E2").Value = -735.46

# Row 3
$ws.Range("A3").Value = "אדיר גז"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1234"
$ws.Range("C3").Value = 45306
$ws.Range("C3").NumberFormat = "m/d/yy"
$ws.Range("D3").Value = 45306
$ws.Range("D3").NumberFormat = "m/d/yy"
$ws.Range("E3").Value = -289.67

# Row 4
$ws.Range("A4").Value = "מופת מילואים"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "20-521-567890"
$ws.Range("C4").Value = 45308
$ws.Range("C4").NumberFormat = "m/d/yy"
$ws.Range("D4").Value = 45308
$ws.Range("D4").NumberFormat = "m/d/yy"
$ws.Range("E4").Value = 1500

# Row 5
$ws.Range("A5").Value = "ביטוח לאומי"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "10-680-335679"
$ws.Range("C5").Value = 45324
$ws.Range("C5").NumberFormat = "m/d/yy"
$ws.Range("D5").Value = 45324
$ws.Range("D5").NumberFormat = "m/d/yy"
$ws.Range("E5").Value = -1300

# Row 6
$ws.Range("A6").Value = "מכבי שירותי בריאות"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "5678"
$ws.Range("C6").Value = 45293
$ws.Range("C6").NumberFormat = "m/d/yy"
$ws.Range("D6").Value = 45293
$ws.Range("D6").NumberFormat = "m/d/yy"
$ws.Range("E6").Value = -257.49

# Row 7
$ws.Range("A7").Value = "מנוי riseup"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "1234"
$ws.Range("C7").Value = 45293
$ws.Range("C7").NumberFormat = "m/d/yy"
$ws.Range("D7").Value = 45293
$ws.Range("D7").NumberFormat = "m/d/yy"
$ws.Range("E7").Value = -45

# Row 8
$ws.Range("A8").Value = "איילון ביטוח חיים"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "5678"
$ws.Range("C8").Value = 45293
$ws.Range("C8").NumberFormat = "m/d/yy"
$ws.Range("D8").Value = 45293
$ws.Range("D8").NumberFormat = "m/d/yy"
$ws.Range("E8").Value = -56.62

# Row 9
$ws.Range("A9").Value = "תרומה חסדי נעמי"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "20-521-567890"
$ws.Range("C9").Value = 45324
$ws.Range("C9").NumberFormat = "m/d/yy"
$ws.Range("D9").Value = 45324
$ws.Range("D9").NumberFormat = "m/d/yy"
$ws.Range("E9").Value = -350

# Row 10
$ws.Range("A10").Value = "הכשרה חובה אקספרס-צמ"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "5678"
$ws.Range("C10").Value = 45324
$ws.Range("C10").NumberFormat = "m/d/yy"
$ws.Range("D10").Value = 45324
$ws.Range("D10").NumberFormat = "m/d/yy"
$ws.Range("E10").Value = -144

# Row 11
$ws.Range("A11").Value = "רמי לוי"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1234"
$ws.Range("C11").Value = 45333
$ws.Range("C11").NumberFormat = "m/d/yy"
$ws.Range("D11").Value = 45333
$ws.Range("D11").NumberFormat = "m/d/yy"
$ws.Range("E11").Value = -327.48

# Row 12
$ws.Range("A12").Value = "משכורת אינטל"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "20-521-567890"
$ws.Range("C12").Value = 45324
$ws.Range("C12").NumberFormat = "m/d/yy"
$ws.Range("D12").Value = 45324
$ws.Range("D12").NumberFormat = "m/d/yy"
$ws.Range("E12").Value = 15700

# Row 13
$ws.Range("A13").Value = "מנורה מבטחים פנסיה"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "10-680-335679"
$ws.Range("C13").Value = 45324
$ws.Range("C13").NumberFormat = "m/d/yy"
$ws.Range("D13").Value = 45324
$ws.Range("D13").NumberFormat = "m/d/yy"
$ws.Range("E13").Value = -439.47

# Row 14
$ws.Range("A14").Value = "מגדל חיים/בריאות"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "1234"
$ws.Range("C14").Value = 45324
$ws.Range("C14").NumberFormat = "m/d/yy"
$ws.Range("D14").Value = 45324
$ws.Range("D14").NumberFormat = "m/d/yy"
$ws.Range("E14").Value = -11.97

# Row 15
$ws.Range("A15").Value = "מגדל חיים/בריאות"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "1234"
$ws.Range("C15").Value = 45333
$ws.Range("C15").NumberFormat = "m/d/yy"
$ws.Range("D15").Value = 45333
$ws.Range("D15").NumberFormat = "m/d/yy"
$ws.Range("E15").Value = -23.35

# Row 16
$ws.Range("A16").Value = "ביטוח ישיר - רכב"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "5678"
$ws.Range("C16").Value = 45324
$ws.Range("C16").NumberFormat = "m/d/yy"
$ws.Range("D16").Value = 45324
$ws.Range("D16").NumberFormat = "m/d/yy"
$ws.Range("E16").Value = -190.25

# Row 17
$ws.Range("A17").Value = "קצבת ילדים"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "20-521-567890"
$ws.Range("C17").Value = 45332
$ws.Range("C17").NumberFormat = "m/d/yy"
$ws.Range("D17").Value = 45332
$ws.Range("D17").NumberFormat = "m/d/yy"
$ws.Range("E17").Value = 550

# Row 18
$ws.Range("A18").Value = "אדיר גז"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "5678"
$ws.Range("C18").Value = 45636
$ws.Range("C18").NumberFormat = "m/d/yy"
$ws.Range("D18").Value = 45636
$ws.Range("D18").NumberFormat = "m/d/yy"
$ws.Range("E18").Value = -323.4

# Row 19
$ws.Range("A19").Value = "צילומינציה"
$ws.Range("B19").Value = 5670
$ws.Range("B19").NumberFormat = "@"
$ws.Range("C19").Value = 45332
$ws.Range("C19").NumberFormat = "m/d/yy"
$ws.Range("D19").Value = 45332
$ws.Range("D19").NumberFormat = "m/d/yy"
$ws.Range("E19").Value = -170

# Row 20
$ws.Range("A20").Value = "צילומינציה"
$ws.Range("B20").Value = 5670
$ws.Range("B20").NumberFormat = "@"
$ws.Range("C20").Value = 45337
$ws.Range("C20").NumberFormat = "m/d/yy"
$ws.Range("D20").Value = 45337
$ws.Range("D20").NumberFormat = "m/d/yy"
$ws.Range("E20").Value = -300


# Column widths
$ws.Columns.Item(1).ColumnWidth = 27.5
$ws.Columns.Item(2).ColumnWidth = 19.296875
$ws.Columns.Item(3).ColumnWidth = 17
$ws.Columns.Item(4).ColumnWidth = 23.3984375

# Row 19-20 B column uses General format + centered (style s=3)
$ws.Range("B19").HorizontalAlignment = -4108
$ws.Range("B20").HorizontalAlignment = -4108

# Selection
$ws.Range("G17").Select()

# Page setup
$ws.PageSetup.Orientation = 1

Write-Host "done"
